# Refresh the cryptos price/volume table (D = Price, E = Volume(1h))
# for the rows whose figures moved in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.859.00"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.638.20"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  +0.60%  "

$ws.Range("D5").Value = "'215.37"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").Value = "'28.75"
$ws.Range("E8").Value = "  -2.72%  "

$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("D11").Value = "'0.0901"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").Value = "1.873.72"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").Value = "1.648.89"
$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("E14").Value = "  +4.15%  "

$ws.Range("E15").Value = "  +7.32%  "

$ws.Range("D16").Value = "'3.87"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").Value = "29.872.39"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "'64.58"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").Value = "'240.34"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").Value = "'9.90"
$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  +2.79%  "

$ws.Range("D25").Value = "'157.70"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").Value = "'15.53"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("E31").Value = "  -0.53%  "

$ws.Range("E32").Value = "  +1.69%  "

$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").Value = "1.424.61"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  +3.61%  "

$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("E37").Value = "  -4.70%  "

$ws.Range("E38").Value = "  +2.52%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "'76.69"
$ws.Range("E40").Value = "  +11.11%  "

$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("E43").Value = "  -1.61%  "

$ws.Range("D44").Value = "'1.96"
$ws.Range("E44").Value = "  -1.12%  "

$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("E46").Value = "  -1.51%  "

$ws.Range("D47").Value = "'5.38"
$ws.Range("E47").Value = "  -0.79%  "

$ws.Range("D48").Value = "1.781.40"
$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").Value = "'48.84"
$ws.Range("E49").Value = "  -9.67%  "

$ws.Range("D50").Value = "'93.38"
$ws.Range("E50").Value = "  +5.91%  "

$ws.Range("E51").Value = "  -0.07%  "
